$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from G1 (antecedents_length) into H1 so the new
# column picks up the same bold font, border and centered alignment.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "consequents_length"

# Fill in the consequents_length values (every rule in this dataset has a
# single-item consequent, so the value is always 1) for rows 2-34.
$lastRow = 34
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}
